# Update cryptocurrency price/volume figures per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '59.632.31'
$ws.Range('E2').Value = '  +0.44%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.650.46'
$ws.Range('E3').Value = '  +1.75%  '
$ws.Range('E4').Value = '  -0.05%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '537.87'
$ws.Range('E5').Value = '  -0.31%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '146.27'
$ws.Range('E6').Value = '  +3.74%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('E8').Value = '  +0.99%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '6.74'
$ws.Range('E9').Value = '  +4.58%  '
$ws.Range('E10').Value = '  -0.04%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.338'
$ws.Range('E11').Value = '  +1.42%  '
$ws.Range('E12').Value = '  -0.19%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '3.114.42'
$ws.Range('E13').Value = '  +1.57%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '59.533.86'
$ws.Range('E14').Value = '  +0.43%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '21.31'
$ws.Range('E15').Value = '  +3.63%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '2.644.95'
$ws.Range('E16').Value = '  +0.59%  '
$ws.Range('E17').Value = '  +1.06%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '340.29'
$ws.Range('E18').Value = '  -0.92%  '
$ws.Range('E19').Value = '  +1.55%  '
$ws.Range('E20').Value = '  +2.26%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '6.24'
$ws.Range('E21').Value = '  -2.81%  '
$ws.Range('E22').Value = '  -0.05%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '66.65'
$ws.Range('E23').Value = '  -1.05%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '0.418'
$ws.Range('E24').Value = '  +2.48%  '
$ws.Range('E25').Value = '  -0.83%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.999'
$ws.Range('E26').Value = '  -0.07%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '7.29'
$ws.Range('E27').Value = '  +1.50%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '0.0₃0747'
$ws.Range('E28').Value = '  +1.32%  '
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('E30').Value = '  -0.87%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '5.85'
$ws.Range('E31').Value = '  +0.24%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '18.92'
$ws.Range('E32').Value = '  +0.86%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '150.91'
$ws.Range('E33').Value = '  +0.97%  '
$ws.Range('E34').Value = '  +0.83%  '
$ws.Range('E35').Value = '  +2.26%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.843'
$ws.Range('E36').Value = '  +3.31%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.839'
$ws.Range('E37').Value = '  +0.81%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '1.45'
$ws.Range('E38').Value = '  -0.92%  '
$ws.Range('E39').Value = '  +1.72%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '286.86'
$ws.Range('E40').Value = '  +4.85%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.998'
$ws.Range('E41').Value = '  -0.08%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.607'
$ws.Range('E42').Value = '  +2.07%  '
$ws.Range('E43').Value = '  -0.07%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.0540'
$ws.Range('E44').Value = '  +3.29%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '19.30'
$ws.Range('E45').Value = '  +4.01%  '
$ws.Range('E46').Value = '  -1.08%  '
$ws.Range('E47').Value = '  +1.89%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.967.66'
$ws.Range('E48').Value = '  +1.20%  '
$ws.Range('E49').Value = '  +1.12%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '18.36'
$ws.Range('E50').Value = '  +0.54%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '111.09'
$ws.Range('E51').Value = '  +0.07%  '
